$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers but must remain text
# (matching the original inlineStr cell type). Mark them as Text first,
# so Excel does not silently convert "304.28" -> 304.28 (number).
$textCells = @("D5", "D6", "D7", "D10", "D11", "D16", "D20", "D23", "D24", "D25", "D28", "D31", "D32", "D33", "D36", "D37", "D41", "D44", "D46", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '41.774.13'
$ws.Range("E2").Value = '  +2.87%  '

$ws.Range("D3").Value = '2.265.75'
$ws.Range("E3").Value = '  +1.46%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '304.28'
$ws.Range("E5").Value = '  +0.79%  '

$ws.Range("D6").Value = '92.00'
$ws.Range("E6").Value = '  +2.82%  '

$ws.Range("D7").Value = '0.530'
$ws.Range("E7").Value = '  +2.37%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("E9").Value = '  +0.86%  '

$ws.Range("D10").Value = '32.33'
$ws.Range("E10").Value = '  +2.48%  '

$ws.Range("D11").Value = '53.19'
$ws.Range("E11").Value = '  +0.87%  '

$ws.Range("E12").Value = '  +1.04%  '

$ws.Range("E13").Value = '  +0.45%  '

$ws.Range("E14").Value = '  +1.62%  '

$ws.Range("D15").Value = '2.617.12'
$ws.Range("E15").Value = '  +1.50%  '

$ws.Range("D16").Value = '14.18'
$ws.Range("E16").Value = '  +1.55%  '

$ws.Range("D17").Value = '2.284.55'
$ws.Range("E17").Value = '  +1.91%  '

$ws.Range("E18").Value = '  +2.68%  '

$ws.Range("D19").Value = '41.681.35'
$ws.Range("E19").Value = '  +2.92%  '

$ws.Range("D20").Value = '12.63'
$ws.Range("E20").Value = '  +10.16%  '

$ws.Range("D21").Value = '0.0₃0904'
$ws.Range("E21").Value = '  +0.96%  '

$ws.Range("E22").Value = '  +1.42%  '

$ws.Range("D23").Value = '66.82'

$ws.Range("D24").Value = '239.85'
$ws.Range("E24").Value = '  +1.08%  '

$ws.Range("D25").Value = '2.58'
$ws.Range("E25").Value = '  +1.70%  '

$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("E27").Value = '  +3.95%  '

$ws.Range("D28").Value = '24.02'
$ws.Range("E28").Value = '  +1.06%  '

$ws.Range("E30").Value = '  -4.33%  '

$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = '34.63'
$ws.Range("E31").Value = '  +5.35%  '

$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").Value = '160.42'
$ws.Range("E32").Value = '  +2.13%  '

$ws.Range("D33").Value = '5.26'
$ws.Range("E33").Value = '  +4.25%  '

$ws.Range("E34").Value = '  -0.22%  '

$ws.Range("E35").Value = '  +2.75%  '

$ws.Range("D36").Value = '3.00'
$ws.Range("E36").Value = '  -0.21%  '

$ws.Range("D37").Value = '16.94'
$ws.Range("E37").Value = '  +6.80%  '

$ws.Range("E38").Value = '  +1.78%  '

$ws.Range("E39").Value = '  +1.75%  '

$ws.Range("E40").Value = '  +0.29%  '

$ws.Range("D41").Value = '1.79'
$ws.Range("E41").Value = '  +1.75%  '

$ws.Range("E42").Value = '  +2.46%  '

$ws.Range("D43").Value = '2.036.49'
$ws.Range("E43").Value = '  -2.57%  '

$ws.Range("D44").Value = '19.32'
$ws.Range("E44").Value = '  -1.58%  '

$ws.Range("E45").Value = '  +1.46%  '

$ws.Range("D46").Value = '10.35'
$ws.Range("E46").Value = '  +2.30%  '

$ws.Range("E47").Value = '  +11.80%  '

$ws.Range("E48").Value = '  +1.00%  '

$ws.Range("E49").Value = '  +0.73%  '

$ws.Range("E50").Value = '  +1.77%  '

$ws.Range("D51").Value = '72.69'
$ws.Range("E51").Value = '  +5.15%  '

# Restore the default (Normal) style on the text-forced cells so they
# keep matching the original (un-styled) look of the other data cells.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
